$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D (Price) and E (Volume) to remain text so numeric-
# and percent-looking strings are not auto-converted to numbers.
$ws.Range("D2:E26").NumberFormat = "@"
$ws.Range("D38:E50").NumberFormat = "@"

$ws.Range("D2").Value = "326.45"
$ws.Range("E2").Value = "-0.62%"
$ws.Range("D3").Value = "44.01"
$ws.Range("E3").Value = "0.14%"
$ws.Range("D4").Value = "5.554"
$ws.Range("E4").Value = "-0.28%"
$ws.Range("D5").Value = "0.07997"
$ws.Range("E5").Value = "-1.41%"
$ws.Range("D6").Value = "1.935"
$ws.Range("E6").Value = "1.53%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "4.327"
$ws.Range("E7").Value = "0.81%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "2.555"
$ws.Range("E8").Value = "-7.48%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9449"
$ws.Range("E9").Value = "-0.45%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1166"
$ws.Range("E10").Value = "-1.17%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1843"
$ws.Range("E11").Value = "-2.99%"
$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D12").Value = "12.07"
$ws.Range("E12").Value = "39.97%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "0.09598"
$ws.Range("E13").Value = "-0.69%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.04778"
$ws.Range("E14").Value = "16.29%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.1066"
$ws.Range("E15").Value = "-0.12%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001279"
$ws.Range("E16").Value = "-0.07%"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04213"
$ws.Range("E17").Value = "-2.42%"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "0.006033"
$ws.Range("E18").Value = "2.05%"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "3.374"
$ws.Range("E19").Value = "-5.58%"
$ws.Range("D20").Value = "0.3474"
$ws.Range("E20").Value = "-0.28%"
$ws.Range("D21").Value = "0.1412"
$ws.Range("E21").Value = "3.51%"
$ws.Range("D22").Value = "0.2517"
$ws.Range("E22").Value = "-2.74%"
$ws.Range("D23").Value = "0.001255"
$ws.Range("E23").Value = "1.31%"
$ws.Range("D24").Value = "0.004305"
$ws.Range("E24").Value = "-2.06%"
$ws.Range("D25").Value = "0.0001198"
$ws.Range("E25").Value = "-2.71%"
$ws.Range("D26").Value = "0.0003762"
$ws.Range("E26").Value = "-5.79%"
$ws.Range("D38").Value = "0.02538"
$ws.Range("E38").Value = "-4.98%"
$ws.Range("D39").Value = "0.05426"
$ws.Range("E39").Value = "-0.79%"
$ws.Range("D40").Value = "0.007553"
$ws.Range("E40").Value = "-1.29%"
$ws.Range("D41").Value = "0.1385"
$ws.Range("E41").Value = "-0.83%"
$ws.Range("D42").Value = "0.007499"
$ws.Range("E42").Value = "-33.84%"
$ws.Range("D43").Value = "0.002029"
$ws.Range("E43").Value = "-3.85%"
$ws.Range("D44").Value = "0.008332"
$ws.Range("E44").Value = "-14.61%"
$ws.Range("D45").Value = "0.00007116"
$ws.Range("E45").Value = "1.39%"
$ws.Range("D46").Value = "0.00000000755"
$ws.Range("E46").Value = "0.54%"
$ws.Range("E47").Value = "1.30%"
$ws.Range("D48").Value = "0.003478"
$ws.Range("E48").Value = "0.73%"
$ws.Range("D49").Value = "0.00002114"
$ws.Range("E49").Value = "0.54%"
$ws.Range("D50").Value = "0.0002013"
$ws.Range("E50").Value = "0.54%"
